# Update the "ToEmail" addresses on the Data sheet from the vajrang@ domain
# test accounts to the abc@ domain test accounts, while preserving the
# existing mailto: hyperlinks attached to each cell (Excel keeps the original
# hyperlink target/uid and records the previous display text in the
# hyperlink's "display" attribute when the cell text is edited in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Capture each hyperlink's current (old) display text before the cell values
# change, then re-assign it as the hyperlink's TextToDisplay so Excel keeps
# it recorded as the "display" attribute even after the cell text changes.
foreach ($h in $ws.Hyperlinks) {
    $h.TextToDisplay = $h.TextToDisplay
}

# New email values for the ToEmail column (A2:A4)
$ws.Range("A2").Value = "abc@outlook.com"
$ws.Range("A3").Value = "abc@gmail.com"
$ws.Range("A4").Value = "abc@yahoo.com"
